$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.22757333333334
$ws.Range("H2").Value = 87.68272
$ws.Range("I2").Value = 0.08948272176993048
$ws.Range("J2").Value = 0.08948272176993047
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 4661.319705626738
$ws.Range("R2").Value = 41951.87735064064
$ws.Range("S2").Value = 0.02669524692126922
$ws.Range("T2").Value = 0.02669524692126921
$ws.Range("G3").Value = 29.22757333333334
$ws.Range("H3").Value = 87.68272
$ws.Range("I3").Value = 0.08948272176993048
$ws.Range("J3").Value = 0.08948272176993047
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 5043.46612535728
$ws.Range("R3").Value = 45391.19512821552
$ws.Range("S3").Value = 0.02888378872467125
$ws.Range("T3").Value = 0.02888378872467125
$ws.Range("G4").Value = 29.22757333333334
$ws.Range("H4").Value = 87.68272
$ws.Range("I4").Value = 0.08948272176993048
$ws.Range("J4").Value = 0.08948272176993047
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 2174.172122470915
$ws.Range("R4").Value = 19567.54910223824
$ws.Range("S4").Value = 0.0124514226279395
$ws.Range("T4").Value = 0.0124514226279395
$ws.Range("G5").Value = 29.22757333333334
$ws.Range("H5").Value = 87.68272
$ws.Range("I5").Value = 0.08948272176993048
$ws.Range("J5").Value = 0.08948272176993047
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 1707.317424165885
$ws.Range("R5").Value = 15365.85681749296
$ws.Range("S5").Value = 0.009777758894348468
$ws.Range("T5").Value = 0.009777758894348468
$ws.Range("G6").Value = 29.22757333333334
$ws.Range("H6").Value = 87.68272
$ws.Range("I6").Value = 0.08948272176993048
$ws.Range("J6").Value = 0.08948272176993047
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 2038.512642862507
$ws.Range("R6").Value = 18346.61378576256
$ws.Range("S6").Value = 0.01167450460170203
$ws.Range("T6").Value = 0.01167450460170203
$ws.Range("I7").Value = 0.4075568457508759
$ws.Range("J7").Value = 0.4075568457508759
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 21230.38636605288
$ws.Range("R7").Value = 191073.4772944759
$ws.Range("S7").Value = 0.1215858259178398
$ws.Range("T7").Value = 0.1215858259178398
$ws.Range("I8").Value = 0.4075568457508759
$ws.Range("J8").Value = 0.4075568457508759
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("Q8").Value = 22970.90550047092
$ws.Range("R8").Value = 206738.1495042383
$ws.Range("S8").Value = 0.1315537300734798
$ws.Range("T8").Value = 0.1315537300734798
$ws.Range("I9").Value = 0.4075568457508759
$ws.Range("J9").Value = 0.4075568457508759
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 9902.456192961887
$ws.Range("R9").Value = 89122.10573665697
$ws.Range("S9").Value = 0.05671108825233991
$ws.Range("T9").Value = 0.05671108825233991
$ws.Range("I10").Value = 0.4075568457508759
$ws.Range("J10").Value = 0.4075568457508759
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 7776.125829940754
$ws.Range("R10").Value = 69985.13246946679
$ws.Range("S10").Value = 0.04453365403590505
$ws.Range("T10").Value = 0.04453365403590505
$ws.Range("I11").Value = 0.4075568457508759
$ws.Range("J11").Value = 0.4075568457508759
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 9284.583283959833
$ws.Range("R11").Value = 83561.2495556385
$ws.Range("S11").Value = 0.05317254747131125
$ws.Range("T11").Value = 0.05317254747131126
$ws.Range("G12").Value = 49.29039633333334
$ws.Range("H12").Value = 147.871189
$ws.Range("I12").Value = 0.1509067745968169
$ws.Range("J12").Value = 0.1509067745968168
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 7861.011692841598
$ws.Range("R12").Value = 70749.10523557437
$ws.Range("S12").Value = 0.04501979298653907
$ws.Range("T12").Value = 0.04501979298653906
$ws.Range("G13").Value = 49.29039633333334
$ws.Range("H13").Value = 147.871189
$ws.Range("I13").Value = 0.1509067745968169
$ws.Range("J13").Value = 0.1509067745968168
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 8505.476707814312
$ws.Range("R13").Value = 76549.2903703288
$ws.Range("S13").Value = 0.04871062601093958
$ws.Range("T13").Value = 0.04871062601093958
$ws.Range("G14").Value = 49.29039633333334
$ws.Range("H14").Value = 147.871189
$ws.Range("I14").Value = 0.1509067745968169
$ws.Range("J14").Value = 0.1509067745968168
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 3666.599494637346
$ws.Range("R14").Value = 32999.39545173612
$ws.Range("S14").Value = 0.02099851223519205
$ws.Range("T14").Value = 0.02099851223519205
$ws.Range("G15").Value = 49.29039633333334
$ws.Range("H15").Value = 147.871189
$ws.Range("I15").Value = 0.1509067745968169
$ws.Range("J15").Value = 0.1509067745968168
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 2879.279492148814
$ws.Range("R15").Value = 25913.51542933933
$ws.Range("S15").Value = 0.01648955271303894
$ws.Range("T15").Value = 0.01648955271303893
$ws.Range("G16").Value = 49.29039633333334
$ws.Range("H16").Value = 147.871189
$ws.Range("I16").Value = 0.1509067745968169
$ws.Range("J16").Value = 0.1509067745968168
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 3437.818629390275
$ws.Range("R16").Value = 30940.36766451248
$ws.Range("S16").Value = 0.01968829065110721
$ws.Range("T16").Value = 0.01968829065110721
$ws.Range("G17").Value = 33.14535033333333
$ws.Range("H17").Value = 99.43605099999999
$ws.Range("I17").Value = 0.1014773319706963
$ws.Range("J17").Value = 0.1014773319706963
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 5286.141031847613
$ws.Range("R17").Value = 47575.26928662851
$ws.Range("S17").Value = 0.03027358109238535
$ws.Range("T17").Value = 0.03027358109238534
$ws.Range("G18").Value = 33.14535033333333
$ws.Range("H18").Value = 99.43605099999999
$ws.Range("I18").Value = 0.1014773319706963
$ws.Range("J18").Value = 0.1014773319706963
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 5719.511836058448
$ws.Range("R18").Value = 51475.60652452603
$ws.Range("S18").Value = 0.03275548350575387
$ws.Range("T18").Value = 0.03275548350575387
$ws.Range("G19").Value = 33.14535033333333
$ws.Range("H19").Value = 99.43605099999999
$ws.Range("I19").Value = 0.1014773319706963
$ws.Range("J19").Value = 0.1014773319706963
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 2465.606564814552
$ws.Range("R19").Value = 22190.45908333097
$ws.Range("S19").Value = 0.01412045948682187
$ws.Range("T19").Value = 0.01412045948682187
$ws.Range("G20").Value = 33.14535033333333
$ws.Range("H20").Value = 99.43605099999999
$ws.Range("I20").Value = 0.1014773319706963
$ws.Range("J20").Value = 0.1014773319706963
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 1936.172856664888
$ws.Range("R20").Value = 17425.55570998399
$ws.Range("S20").Value = 0.011088407522989
$ws.Range("T20").Value = 0.011088407522989
$ws.Range("G21").Value = 33.14535033333333
$ws.Range("H21").Value = 99.43605099999999
$ws.Range("I21").Value = 0.1014773319706963
$ws.Range("J21").Value = 0.1014773319706963
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 2311.762763744339
$ws.Range("R21").Value = 20805.86487369905
$ws.Range("S21").Value = 0.01323940036274625
$ws.Range("T21").Value = 0.01323940036274625
$ws.Range("G22").Value = 81.84527466666667
$ws.Range("H22").Value = 245.535824
$ws.Range("I22").Value = 0.2505763259116804
$ws.Range("J22").Value = 0.2505763259116804
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 13052.98210238572
$ws.Range("R22").Value = 117476.8389214715
$ws.Range("S22").Value = 0.07475406157219233
$ws.Range("T22").Value = 0.07475406157219233
$ws.Range("G23").Value = 81.84527466666667
$ws.Range("H23").Value = 245.535824
$ws.Range("I23").Value = 0.2505763259116804
$ws.Range("J23").Value = 0.2505763259116804
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 14123.09758303217
$ws.Range("R23").Value = 127107.8782472896
$ws.Range("S23").Value = 0.08088258284818337
$ws.Range("T23").Value = 0.08088258284818337
$ws.Range("G24").Value = 81.84527466666667
$ws.Range("H24").Value = 245.535824
$ws.Range("I24").Value = 0.2505763259116804
$ws.Range("J24").Value = 0.2505763259116804
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 6088.28220211149
$ws.Range("R24").Value = 54794.5398190034
$ws.Range("S24").Value = 0.0348674210257548
$ws.Range("T24").Value = 0.0348674210257548
$ws.Range("G25").Value = 81.84527466666667
$ws.Range("H25").Value = 245.535824
$ws.Range("I25").Value = 0.2505763259116804
$ws.Range("J25").Value = 0.2505763259116804
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 4780.96015465907
$ws.Range("R25").Value = 43028.64139193163
$ws.Range("S25").Value = 0.0273804244097033
$ws.Range("T25").Value = 0.0273804244097033
$ws.Range("G26").Value = 81.84527466666667
$ws.Range("H26").Value = 245.535824
$ws.Range("I26").Value = 0.2505763259116804
$ws.Range("J26").Value = 0.2505763259116804
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 5708.398205480662
$ws.Range("R26").Value = 51375.58384932595
$ws.Range("S26").Value = 0.03269183605584657
$ws.Range("T26").Value = 0.03269183605584658
